# Applies the Zeromus_Profits Leve-profit recalculation update across all
# job tables (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR): refreshed market-board
# price/profit figures (columns H-N) for the affected Leve rows.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 12 (Leve Item ID 5515)
$ws.Range("H12").Value = 1057.5
$ws.Range("I12").Value = 132.66667
$ws.Range("J12").Value = 1365.7778
$ws.Range("K12").Value = 132.66667
$ws.Range("L12").Value = 1365.7778
$ws.Range("M12").Value = 37.33332999999999
$ws.Range("N12").Value = -1705.7778
# Row 100 (Leve Item ID 19906)
$ws.Range("H100").Value = 50004350
$ws.Range("I100").Value = 100001750
$ws.Range("J100").Value = 6950
$ws.Range("K100").Value = 100001750
$ws.Range("L100").Value = 6950
$ws.Range("M100").Value = -100001209
$ws.Range("N100").Value = -8032

$ws = $wb.Worksheets.Item("ARM")
# Row 24 (Leve Item ID 18363)
$ws.Range("H24").Value = 30000
$ws.Range("J24").Value = 30000
$ws.Range("L24").Value = 30000
$ws.Range("N24").Value = -30748
# Row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 11489.569
$ws.Range("I32").Value = 3034.98
$ws.Range("J32").Value = 30704.545
$ws.Range("K32").Value = 3034.98
$ws.Range("L32").Value = 30704.545
$ws.Range("M32").Value = -2747.98
$ws.Range("N32").Value = -31278.545
# Row 45 (Leve Item ID 27714)
$ws.Range("H45").Value = 1209.3334
$ws.Range("I45").Value = 1135.5
$ws.Range("K45").Value = 1135.5
$ws.Range("M45").Value = -758.5
# Row 86 (Leve Item ID 10702)
$ws.Range("H86").Value = 50000
$ws.Range("I86").Value = 50000
$ws.Range("K86").Value = 50000
$ws.Range("M86").Value = -48814
# Row 89 (Leve Item ID 10702)
$ws.Range("H89").Value = 50000
$ws.Range("I89").Value = 50000
$ws.Range("K89").Value = 150000
$ws.Range("M89").Value = -144072
# Row 94 (Leve Item ID 18055)
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
# Row 97 (Leve Item ID 19941)
$ws.Range("H97").Value = 5954123
$ws.Range("I97").Value = 6668535
$ws.Range("J97").Value = 690
$ws.Range("K97").Value = 6668535
$ws.Range("L97").Value = 690
$ws.Range("M97").Value = -6668039
$ws.Range("N97").Value = -1682
# Row 100 (Leve Item ID 18363)
$ws.Range("H100").Value = 30000
$ws.Range("J100").Value = 30000
$ws.Range("L100").Value = 30000
$ws.Range("N100").Value = -32164
# Row 102 (Leve Item ID 19945)
$ws.Range("H102").Value = 66668916
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 94 (Leve Item ID 19939)
$ws.Range("H94").Value = 9005.791999999999
$ws.Range("I94").Value = 640.7727
$ws.Range("J94").Value = 101021
$ws.Range("K94").Value = 640.7727
$ws.Range("L94").Value = 101021
$ws.Range("M94").Value = -189.7727
$ws.Range("N94").Value = -101923
# Row 99 (Leve Item ID 19943)
$ws.Range("H99").Value = 1720.1666
$ws.Range("I99").Value = 1762
$ws.Range("J99").Value = 1511
$ws.Range("K99").Value = 1762
$ws.Range("L99").Value = 1511
$ws.Range("M99").Value = -264
$ws.Range("N99").Value = -4507
# Row 103 (Leve Item ID 18514)
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
# Row 107 (Leve Item ID 27706)
$ws.Range("H107").Value = 1359.6207
$ws.Range("I107").Value = 1044.3889
$ws.Range("J107").Value = 1875.4546
$ws.Range("K107").Value = 1044.3889
$ws.Range("L107").Value = 1875.4546
$ws.Range("M107").Value = 875.6111000000001
$ws.Range("N107").Value = -5715.4546

$ws = $wb.Worksheets.Item("CRP")
# Row 97 (Leve Item ID 19730)
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("N97").ClearContents()
# Row 103 (Leve Item ID 19558)
$ws.Range("H103").Value = 4999
$ws.Range("I103").Value = 4999
$ws.Range("K103").Value = 4999
$ws.Range("M103").Value = -3827
# Row 122 (Leve Item ID 36196)
$ws.Range("H122").Value = 6061490
$ws.Range("I122").Value = 13334038
$ws.Range("J122").Value = 1033.3334
$ws.Range("K122").Value = 40002114
$ws.Range("L122").Value = 3100.0002
$ws.Range("M122").Value = -39999664
$ws.Range("N122").Value = -8000.0002
# Row 134 (Leve Item ID 44020)
$ws.Range("H134").Value = 7474.778
$ws.Range("I134").Value = 8324.714
$ws.Range("J134").Value = 4500
$ws.Range("K134").Value = 24974.142
$ws.Range("L134").Value = 13500
$ws.Range("M134").Value = -22439.142
$ws.Range("N134").Value = -18570
# Row 141 (Leve Item ID 43345)
$ws.Range("H141").Value = 46080.418
$ws.Range("J141").Value = 48905.91
$ws.Range("L141").Value = 48905.91
$ws.Range("N141").Value = -59265.91

$ws = $wb.Worksheets.Item("CUL")
# Row 3 (Leve Item ID 44094)
$ws.Range("H3").Value = 8545.444
$ws.Range("J3").Value = 8988.625
$ws.Range("L3").Value = 26965.875
$ws.Range("N3").Value = -27189.875
# Row 132 (Leve Item ID 43972)
$ws.Range("H132").Value = 1030.1538
$ws.Range("I132").Value = 913.1429000000001
$ws.Range("J132").Value = 1166.6666
$ws.Range("K132").Value = 8218.286100000001
$ws.Range("L132").Value = 10499.9994
$ws.Range("M132").Value = -5688.286100000001
$ws.Range("N132").Value = -15559.9994
# Row 134 (Leve Item ID 44074)
$ws.Range("H134").Value = 1587.5
$ws.Range("I134").Value = 1587.5
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 4762.5
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = 307.5
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 102 (Leve Item ID 36169)
$ws.Range("H102").Value = 2187
$ws.Range("I102").Value = 2187
$ws.Range("K102").Value = 2187
$ws.Range("M102").Value = -565

$ws = $wb.Worksheets.Item("LTW")
# Row 22 (Leve Item ID 5277)
$ws.Range("H22").Value = 1852277.5
$ws.Range("I22").Value = 3030629.2
$ws.Range("J22").Value = 581.7143
$ws.Range("K22").Value = 3030629.2
$ws.Range("L22").Value = 581.7143
$ws.Range("M22").Value = -3030334.2
$ws.Range("N22").Value = -1171.7143
# Row 27 (Leve Item ID 5277)
$ws.Range("H27").Value = 1852277.5
$ws.Range("I27").Value = 3030629.2
$ws.Range("J27").Value = 581.7143
$ws.Range("K27").Value = 3030629.2
$ws.Range("L27").Value = 581.7143
$ws.Range("M27").Value = -3030522.2
$ws.Range("N27").Value = -795.7143
# Row 46 (Leve Item ID 5282)
$ws.Range("H46").Value = 1976.2222
$ws.Range("I46").Value = 1100.5
$ws.Range("J46").Value = 2226.4285
$ws.Range("K46").Value = 1100.5
$ws.Range("L46").Value = 2226.4285
$ws.Range("M46").Value = -912.5
$ws.Range("N46").Value = -2602.4285
# Row 132 (Leve Item ID 44058)
$ws.Range("H132").Value = 1282057.6
$ws.Range("I132").Value = 2024670.1
$ws.Range("J132").Value = 3113.9443
$ws.Range("K132").Value = 6074010.300000001
$ws.Range("L132").Value = 9341.832900000001
$ws.Range("M132").Value = -6071480.300000001
$ws.Range("N132").Value = -14401.8329

$ws = $wb.Worksheets.Item("WVR")
# Row 98 (Leve Item ID 18374)
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
# Row 100 (Leve Item ID 19981)
$ws.Range("H100").Value = 27637.957
$ws.Range("I100").Value = 36168.707
$ws.Range("J100").Value = 3467.5
$ws.Range("K100").Value = 72337.414
$ws.Range("L100").Value = 6935
$ws.Range("M100").Value = -71796.414
$ws.Range("N100").Value = -8017
# Row 103 (Leve Item ID 18548)
$ws.Range("H103").Value = 34000
$ws.Range("J103").Value = 34000
$ws.Range("L103").Value = 34000
$ws.Range("N103").Value = -36344
# Row 105 (Leve Item ID 18710)
$ws.Range("H105").Value = 39923
$ws.Range("J105").Value = 39923
$ws.Range("L105").Value = 39923
$ws.Range("N105").Value = -46911

Write-Output "Applied all Zeromus_Profits numeric updates."